$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-07-05T16:08:03"
$ws.Range("T4").Value = 61.43
$ws.Range("U4").Value = 80.88
$ws.Range("V4").Value = 123.86
$ws.Range("W4").Value = 64.05
$ws.Range("X4").Value = 46.89
$ws.Range("Y4").Value = 47.98
$ws.Range("Z4").Value = 35.43
$ws.Range("T6").Value = -3.38
$ws.Range("U6").Value = -3.72
$ws.Range("V6").Value = -6.44
$ws.Range("W6").Value = -3.07
$ws.Range("Y6").Value = -0.86
$ws.Range("Z6").Value = -0.78
$ws.Range("T9").Value = 61.02
$ws.Range("U9").Value = 79.81
$ws.Range("V9").Value = 120.88
$ws.Range("W9").Value = 63.99
$ws.Range("X9").Value = 48.65
$ws.Range("Y9").Value = 49.94
$ws.Range("Z9").Value = 37.06
$ws.Range("T11").Value = -3.78
$ws.Range("U11").Value = -4.79
$ws.Range("V11").Value = -9.43
$ws.Range("W11").Value = -3.14
$ws.Range("X11").Value = 0.49
$ws.Range("Y11").Value = 1.1
$ws.Range("Z11").Value = 0.85
$ws.Range("T14").Value = 61.02
$ws.Range("U14").Value = 79.81
$ws.Range("V14").Value = 120.88
$ws.Range("W14").Value = 63.99
$ws.Range("X14").Value = 48.65
$ws.Range("Y14").Value = 49.94
$ws.Range("Z14").Value = 37.06
$ws.Range("T16").Value = -3.78
$ws.Range("U16").Value = -4.79
$ws.Range("V16").Value = -9.43
$ws.Range("W16").Value = -3.14
$ws.Range("X16").Value = 0.49
$ws.Range("Y16").Value = 1.1
$ws.Range("Z16").Value = 0.85
$ws.Range("T19").Value = 61.19
$ws.Range("U19").Value = 80.8
$ws.Range("V19").Value = 99.34999999999999
$ws.Range("W19").Value = 48.92
$ws.Range("X19").Value = 46.76
$ws.Range("Y19").Value = 47.89
$ws.Range("Z19").Value = 35.53
$ws.Range("V20").Value = -23.58
$ws.Range("W20").Value = -14.7
$ws.Range("T21").Value = -3.61
$ws.Range("U21").Value = -3.8
$ws.Range("V21").Value = -7.38
$ws.Range("W21").Value = -3.5
$ws.Range("X21").Value = -1.4
$ws.Range("Y21").Value = -0.96
$ws.Range("Z21").Value = -0.68
$ws.Range("T24").Value = 61.19
$ws.Range("U24").Value = 80.8
$ws.Range("V24").Value = 122.93
$ws.Range("W24").Value = 63.62
$ws.Range("X24").Value = 46.76
$ws.Range("Y24").Value = 47.89
$ws.Range("Z24").Value = 35.53
$ws.Range("T26").Value = -3.61
$ws.Range("U26").Value = -3.8
$ws.Range("V26").Value = -7.38
$ws.Range("W26").Value = -3.5
$ws.Range("X26").Value = -1.4
$ws.Range("Y26").Value = -0.96
$ws.Range("Z26").Value = -0.68
$ws.Range("T29").Value = 60.91
$ws.Range("U29").Value = 80.8
$ws.Range("V29").Value = 98.31
$ws.Range("W29").Value = 48.32
$ws.Range("X29").Value = 46.44
$ws.Range("Y29").Value = 47.56
$ws.Range("Z29").Value = 35.53
$ws.Range("V30").Value = -23.58
$ws.Range("W30").Value = -14.7
$ws.Range("T31").Value = -3.9
$ws.Range("U31").Value = -3.8
$ws.Range("V31").Value = -8.41
$ws.Range("W31").Value = -4.1
$ws.Range("X31").Value = -1.72
$ws.Range("Y31").Value = -1.28
$ws.Range("Z31").Value = -0.68
$ws.Range("T34").Value = 60.74
$ws.Range("U34").Value = 79.20999999999999
$ws.Range("V34").Value = 119.65
$ws.Range("W34").Value = 64.17
$ws.Range("X34").Value = 49.75
$ws.Range("Y34").Value = 50.88
$ws.Range("Z34").Value = 37.68
$ws.Range("T36").Value = -4.07
$ws.Range("U36").Value = -5.39
$ws.Range("V36").Value = -10.65
$ws.Range("W36").Value = -2.95
$ws.Range("X36").Value = 1.59
$ws.Range("Y36").Value = 2.04
$ws.Range("Z36").Value = 1.47
$ws.Range("T39").Value = 61.43
$ws.Range("U39").Value = 80.88
$ws.Range("V39").Value = 123.86
$ws.Range("W39").Value = 64.05
$ws.Range("X39").Value = 46.89
$ws.Range("Y39").Value = 47.98
$ws.Range("Z39").Value = 35.43
$ws.Range("T41").Value = -3.38
$ws.Range("U41").Value = -3.72
$ws.Range("V41").Value = -6.44
$ws.Range("W41").Value = -3.07
$ws.Range("Y41").Value = -0.86
$ws.Range("Z41").Value = -0.78
$ws.Range("T44").Value = 62.67
$ws.Range("U44").Value = 82.45999999999999
$ws.Range("V44").Value = 129.14
$ws.Range("X44").Value = 48.21
$ws.Range("Y44").Value = 49.04
$ws.Range("Z44").Value = 35.64
$ws.Range("T46").Value = -2.13
$ws.Range("U46").Value = -2.14
$ws.Range("V46").Value = -1.16
$ws.Range("W46").Value = -0.4
$ws.Range("Y46").Value = 0.2
$ws.Range("Z46").Value = -0.57
$ws.Range("T49").Value = 61.31
$ws.Range("U49").Value = 79.73999999999999
$ws.Range("V49").Value = 125.53
$ws.Range("W49").Value = 64.48
$ws.Range("X49").Value = 48.94
$ws.Range("Y49").Value = 49.34
$ws.Range("Z49").Value = 35.85
$ws.Range("T51").Value = -3.49
$ws.Range("U51").Value = -4.86
$ws.Range("V51").Value = -4.77
$ws.Range("W51").Value = -2.64
$ws.Range("X51").Value = 0.78
$ws.Range("Y51").Value = 0.49
$ws.Range("Z51").Value = -0.36
$ws.Range("T54").Value = 60.17
$ws.Range("U54").Value = 78.62
$ws.Range("V54").Value = 122.47
$ws.Range("W54").Value = 65.04000000000001
$ws.Range("X54").Value = 46.94
$ws.Range("Y54").Value = 48.03
$ws.Range("Z54").Value = 35.57
$ws.Range("T56").Value = -4.63
$ws.Range("U56").Value = -5.98
$ws.Range("V56").Value = -7.84
$ws.Range("W56").Value = -2.08
$ws.Range("Y56").Value = -0.82
$ws.Range("Z56").Value = -0.64
$ws.Range("U59").Value = 84.68000000000001
$ws.Range("V59").Value = 133.78
$ws.Range("W59").Value = 69.2
$ws.Range("X59").Value = 49.6
$ws.Range("Y59").Value = 50.46
$ws.Range("Z59").Value = 36.43
$ws.Range("T61").Value = -0.06
$ws.Range("U61").Value = 0.08
$ws.Range("V61").Value = 3.48
$ws.Range("W61").Value = 2.08
$ws.Range("X61").Value = 1.44
$ws.Range("Y61").Value = 1.61
$ws.Range("Z61").Value = 0.22
$ws.Range("T64").Value = 65.98999999999999
$ws.Range("U64").Value = 86.42
$ws.Range("V64").Value = 136.59
$ws.Range("W64").Value = 70.58
$ws.Range("X64").Value = 50.48
$ws.Range("Y64").Value = 51.25
$ws.Range("Z64").Value = 36.87
$ws.Range("T66").Value = 1.19
$ws.Range("U66").Value = 1.81
$ws.Range("V66").Value = 6.28
$ws.Range("W66").Value = 3.46
$ws.Range("X66").Value = 2.32
$ws.Range("Y66").Value = 2.41
$ws.Range("Z66").Value = 0.66
$ws.Range("U69").Value = 84.77
$ws.Range("V69").Value = 136.73
$ws.Range("W69").Value = 70.95
$ws.Range("X69").Value = 50.69
$ws.Range("Y69").Value = 51.69
$ws.Range("Z69").Value = 36.76
$ws.Range("T71").Value = 0.19
$ws.Range("U71").Value = 0.17
$ws.Range("V71").Value = 6.43
$ws.Range("W71").Value = 3.83
$ws.Range("X71").Value = 2.53
$ws.Range("Y71").Value = 2.84
$ws.Range("Z71").Value = 0.55
$ws.Range("U74").Value = 83.84999999999999
$ws.Range("V74").Value = 132.29
$ws.Range("W74").Value = 68.77
$ws.Range("X74").Value = 49.14
$ws.Range("Y74").Value = 50.05
$ws.Range("Z74").Value = 36.03
$ws.Range("T76").Value = -0.45
$ws.Range("U76").Value = -0.75
$ws.Range("V76").Value = 1.98
$ws.Range("X76").Value = 0.98
$ws.Range("Y76").Value = 1.2
$ws.Range("Z76").Value = -0.18
$ws.Range("T79").Value = 64.58
$ws.Range("U79").Value = 84.27
$ws.Range("V79").Value = 132.77
$ws.Range("W79").Value = 69.05
$ws.Range("X79").Value = 49.36
$ws.Range("Y79").Value = 50.3
$ws.Range("Z79").Value = 36.29
$ws.Range("T81").Value = -0.22
$ws.Range("U81").Value = -0.33
$ws.Range("V81").Value = 2.46
$ws.Range("W81").Value = 1.92
$ws.Range("X81").Value = 1.2
$ws.Range("Y81").Value = 1.46
$ws.Range("Z81").Value = 0.08
$ws.Range("T84").Value = 59.24
$ws.Range("U84").Value = 77.61
$ws.Range("V84").Value = 120.88
$ws.Range("W84").Value = 65.87
$ws.Range("X84").Value = 47.4
$ws.Range("Y84").Value = 48.41
$ws.Range("Z84").Value = 35.57
$ws.Range("T86").Value = -5.57
$ws.Range("U86").Value = -6.99
$ws.Range("V86").Value = -9.43
$ws.Range("W86").Value = -1.25
$ws.Range("Y86").Value = -0.44
$ws.Range("Z86").Value = -0.64
$ws.Range("T89").Value = 60.91
$ws.Range("U89").Value = 80.8
$ws.Range("V89").Value = 98.31
$ws.Range("W89").Value = 48.32
$ws.Range("X89").Value = 46.44
$ws.Range("Y89").Value = 47.56
$ws.Range("Z89").Value = 35.53
$ws.Range("V90").Value = -23.58
$ws.Range("W90").Value = -14.7
$ws.Range("T91").Value = -3.9
$ws.Range("U91").Value = -3.8
$ws.Range("V91").Value = -8.41
$ws.Range("W91").Value = -4.1
$ws.Range("X91").Value = -1.72
$ws.Range("Y91").Value = -1.28
$ws.Range("Z91").Value = -0.68
